$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S (19) values for rows 4-34. Column R (18) holds the matching
# 2021 figures; S holds the new 2022 figures added alongside them.
$values = [ordered]@{
    4  = 2022
    5  = 10.071559327675153
    6  = 10.551906067345987
    7  = 9.5619606820956751
    8  = 8.2747510251903922
    9  = 7.6325088339222615
    10 = 8.9652028567087072
    11 = 8.5830821067565175
    12 = 10.275380189066995
    13 = 6.7661984261234096
    14 = 9.0818473806623103
    15 = 9.0186815546489161
    16 = 9.149130832570906
    17 = 8.0270384452893957
    18 = 8.8235294117647065
    19 = 7.2217502124044177
    20 = 3.4213262670647033
    21 = 3.4802022457154114
    22 = 3.3598464070213931
    23 = 12.808072967203572
    24 = 14.988470407378941
    25 = 10.584084672677381
    26 = 7.1442946266854497
    27 = 7.5305623471882646
    28 = 6.7350533623458704
    29 = 16.241806263656226
    30 = 15.614010409340272
    31 = 16.915873735085334
    32 = 22.910065805508165
    33 = 24.889729048519218
    34 = 20.801878879382652
}

foreach ($row in $values.Keys) {
    $rCell = $ws.Cells.Item($row, 18)
    $sCell = $ws.Cells.Item($row, 19)

    # Copy column R's formatting (style) into the new column S cell so the
    # cell picks up the same cell style index rather than creating a new one.
    $rCell.Copy() | Out-Null
    $sCell.PasteSpecial(-4122) | Out-Null

    $sCell.Value = $values[$row]
}

$ws.Range("A1").Select()
